$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '27.678.31'
$ws.Cells.Item(2, 5).Value = '  +0.51%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.638.86'
$ws.Cells.Item(3, 5).Value = '  -0.47%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '212.69'
$ws.Cells.Item(5, 5).Value = '  +0.04%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -1.27%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '23.05'
$ws.Cells.Item(8, 5).Value = '  -2.21%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.258'
$ws.Cells.Item(9, 5).Value = '  +0.03%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0611'
$ws.Cells.Item(10, 5).Value = '  -0.21%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0895'
$ws.Cells.Item(11, 5).Value = '  +0.37%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.870.03'
$ws.Cells.Item(12, 5).Value = '  -0.49%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.648.23'
$ws.Cells.Item(13, 5).Value = '  +0.28%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '4.05'
$ws.Cells.Item(14, 5).Value = '  -0.12%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.559'
$ws.Cells.Item(15, 5).Value = '  -5.67%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.66'
$ws.Cells.Item(16, 5).Value = '  +0.16%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '27.655.09'
$ws.Cells.Item(17, 5).Value = '  +0.55%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '230.45'
$ws.Cells.Item(18, 5).Value = '  -0.70%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +2.95%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0₃0723'
$ws.Cells.Item(20, 5).Value = '  -0.22%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.06%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.31'
$ws.Cells.Item(22, 5).Value = '  -0.81%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.18'
$ws.Cells.Item(23, 5).Value = '  +4.24%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.02'
$ws.Cells.Item(24, 5).Value = '  -0.12%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '150.62'
$ws.Cells.Item(25, 5).Value = '  +1.52%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '6.96'
$ws.Cells.Item(26, 5).Value = '  -1.27%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -1.77%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '15.63'
$ws.Cells.Item(28, 5).Value = '  -0.01%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'BinanceUSD'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  -0.05%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.21%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.01%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.31'
$ws.Cells.Item(32, 5).Value = '  -0.33%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.456.69'
$ws.Cells.Item(33, 5).Value = '  +2.10%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -1.54%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.87%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.23%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.566'
$ws.Cells.Item(37, 5).Value = '  -0.66%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.881'
$ws.Cells.Item(38, 5).Value = '  -1.03%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.17%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.901'
$ws.Cells.Item(40, 5).Value = '  +10.10%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '70.06'
$ws.Cells.Item(41, 5).Value = '  +7.61%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.06%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.97%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.62'
$ws.Cells.Item(44, 5).Value = '  +1.08%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.03%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -0.47%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.779.81'
$ws.Cells.Item(47, 5).Value = '  -0.53%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.74'
$ws.Cells.Item(48, 5).Value = '  +3.13%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '86.70'
$ws.Cells.Item(49, 5).Value = '  -1.75%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.54%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0991'
$ws.Cells.Item(51, 5).Value = '  -0.17%  '
